$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list-column"

# --- Add the new sheet right after the first one ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "two-row-header"

# --- Populate header rows 1 & 2 in the exact order the strings were
#     first typed, so the shared-string table comes out in the same
#     order as the authored workbook. ---
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "species"
$ws2.Range("C1").Value = "death"
$ws2.Range("D1").Value = "weight"
$ws2.Range("D2").Value = "(in grams)"
$ws2.Range("B2").Value = "(office supply type)"
$ws2.Range("A2").Value = "(at birth)"
$ws2.Range("C2").Value = "(date is approximate)"

# --- Data row ---
$ws2.Range("A3").Value = "Clippy"
$ws2.Range("B3").Value = "paperclip"

# Copy the existing date-format style from list-column!B4 onto C3 before
# writing the value, so it reuses cellXfs index 1 instead of minting a
# brand-new number format.
$ws1.Range("B4").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws2.Range("C3").Value = (Get-Date -Year 2007 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)

$ws2.Range("D3").Value = 0.9

# --- Selections: list-column keeps A2:A5 selected and is no longer the
#     active tab; two-row-header becomes the active tab with A1:D1 selected. ---
$ws1.Range("A2:A5").Select()
$ws2.Activate()
$ws2.Range("A1:D1").Select()
